# Insert a new weekly record row above row 577 (the first "Zanahoria"
# price record for "Feria Lagunitas de Puerto Montt"), which shifts all
# subsequent records down by one row (577 -> 578, ..., 693 -> 694).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 577 (pushes existing rows 577..693 down to 578..694)
$ws.Rows.Item(577).Insert()

# Populate the newly inserted row with this week's data
$ws.Range("A577").Value = 4
$ws.Range("B577").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C577").Value = "Los Lagos"
$ws.Range("D577").Value = 45275
$ws.Range("E577").Value = 10
$ws.Range("F577").Value = 100114013
$ws.Range("G577").Value = "Zanahoria"
$ws.Range("H577").Value = "Sin especificar"
$ws.Range("I577").Value = "Primera"
$ws.Range("J577").Value = 700
$ws.Range("K577").Value = 8500
$ws.Range("L577").Value = 9000
$ws.Range("M577").Value = 8750
$ws.Range("N577").Value = "`$/saco 20 kilos"
$ws.Range("O577").Value = "Región Metropolitana"
$ws.Range("P577").Value = 438
$ws.Range("Q577").Value = 20
$ws.Range("R577").Value = "Hortaliza"
